$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Book-level view tweak (tab ratio of the sheet-tab / horizontal-scroll
# splitter). Best effort - harmless if the host doesn't persist window
# chrome state.
try {
    $excel.ActiveWindow.TabRatio = 0.236
} catch {}

# --- Row 2: Worldcheck for SubIs ---
# Highlight Title cell (C2) green
$ws.Range("C2").Interior.Color = 5296274
# Assignee changes from "Pradeep" to "Pradeep " (trailing space)
$ws.Range("J2").Value = "Pradeep "
# New "Done" status added in column K
$ws.Range("K2").Value = "Done"

# --- Row 20: User Privilege fluctuating ---
# Highlight Title cell (C20) green
$ws.Range("C20").Interior.Color = 5296274

# --- Row 22: Chronological Findings display ---
# Highlight Title cell (C22) green
$ws.Range("C22").Interior.Color = 5296274
# Status moves from Pending to Done
$ws.Range("K22").Value = "Done"

# --- Row 24: Live Search to DB Search Conversion ---
# New "Pending" status added in column K
$ws.Range("K24").Value = "Pending"

# --- Row 30: Automated Country Specific weblink inclusion ---
# Highlight Title cell (C30) green
$ws.Range("C30").Interior.Color = 5296274
# Status moves from Pending to Done
$ws.Range("K30").Value = "Done"

# --- Row 31: Automated Sponsor Specific weblink inclusion ---
# Highlight Title cell (C31) green
$ws.Range("C31").Interior.Color = 5296274
# Status moves from Pending to Done
$ws.Range("K31").Value = "Done"

# Scroll the view down one row (so row 28 becomes the top visible row) and
# move the active selection down to C31, matching the author's final
# cursor position. Best effort - the scroll position may not be persisted
# by every host.
try {
    $excel.ActiveWindow.ScrollRow = 28
} catch {}
$ws.Range("C31").Select()
